# ---------------------------------------------------------------------------
# Scheduled market-data refresh for the Leve-profit workbook.
#
# Each worksheet (one per crafting class: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
# lists Leves with columns H:N driven by the latest Universalis market-board
# pull:
#   H currentAveragePrice    I currentAveragePriceNQ  J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ
#   M LeveProfitNQ           N LeveProfitHQ
#
# This run refreshes the rows whose market prices moved since the last pull.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11: Gotta Bounce / Rubber
$ws.Range("H11").Value = 37213.43
$ws.Range("I11").Value = 37213.43
$ws.Range("K11").Value = 37213.43
$ws.Range("M11").Value = -37073.43
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 1099.75
$ws.Range("J17").Value = 1298.2
$ws.Range("L17").Value = 3894.6
$ws.Range("N17").Value = -4230.6
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 961.63635
$ws.Range("I28").Value = 985.26666
$ws.Range("K28").Value = 985.26666
$ws.Range("M28").Value = -500.26666
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 1125.5
$ws.Range("I40").Value = 947.875
$ws.Range("J40").Value = 1480.75
$ws.Range("K40").Value = 947.875
$ws.Range("L40").Value = 1480.75
$ws.Range("M40").Value = -772.875
$ws.Range("N40").Value = -1830.75
# Row 42: Eye of the Beholder / Hi-Potion of Dexterity
$ws.Range("H42").Value = 242
$ws.Range("I42").Value = 219.42857
$ws.Range("K42").Value = 658.28571
$ws.Range("M42").Value = -428.28571
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 18524004
$ws.Range("I137").Value = 20835180
$ws.Range("J137").Value = 34601
$ws.Range("K137").Value = 62505540
$ws.Range("L137").Value = 103803
$ws.Range("M137").Value = -62502990
$ws.Range("N137").Value = -108903
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3474.7144
$ws.Range("J138").Value = 3759.8667
$ws.Range("L138").Value = 11279.6001
$ws.Range("N138").Value = -21559.6001

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate / Bronze Plate
$ws.Range("H4").Value = 2.3333333
# Row 6: Don't Hit Me One More Time / Bronze Hoplon
$ws.Range("H6").Value = 4000
$ws.Range("J6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5346
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 885528.4399999999
$ws.Range("I32").Value = 993763.9
$ws.Range("K32").Value = 993763.9
$ws.Range("M32").Value = -993476.9
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4549580.5
$ws.Range("I61").Value = 4753
$ws.Range("K61").Value = 4753
$ws.Range("M61").Value = -4541
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1116.5454
$ws.Range("I110").Value = 1108.2
$ws.Range("K110").Value = 1108.2
$ws.Range("M110").Value = 936.8
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2786.0894
$ws.Range("I132").Value = 2026
$ws.Range("K132").Value = 6078
$ws.Range("M132").Value = -3548
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4549580.5
$ws.Range("I136").Value = 4753
$ws.Range("K136").Value = 14259
$ws.Range("M136").Value = -11709

$ws = $wb.Worksheets.Item("BSM")
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 3306.4707
$ws.Range("I105").Value = 2400.077
$ws.Range("K105").Value = 2400.077
$ws.Range("M105").Value = -653.0770000000002
# Row 114: Halfhearted Effort / Bluespirit Halfheart Saw
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 3883839.5
$ws.Range("I134").Value = 5297.2285
$ws.Range("K134").Value = 15891.6855
$ws.Range("M134").Value = -13356.6855

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 81.56521600000001
$ws.Range("I7").Value = 43
$ws.Range("K7").Value = 43
$ws.Range("M7").Value = 70
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 5015551.5
$ws.Range("I58").Value = 4169280
$ws.Range("J58").Value = 11785723
$ws.Range("K58").Value = 4169280
$ws.Range("L58").Value = 11785723
$ws.Range("M58").Value = -4169077
$ws.Range("N58").Value = -11786129
# Row 108: Just Starting Out / White Oak Fishing Rod
$ws.Range("H108").Value = 21333.334
$ws.Range("J108").Value = 21333.334
$ws.Range("L108").Value = 21333.334
$ws.Range("N108").Value = -29013.334
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 5015551.5
$ws.Range("I136").Value = 4169280
$ws.Range("J136").Value = 11785723
$ws.Range("K136").Value = 12507840
$ws.Range("L136").Value = 35357169
$ws.Range("M136").Value = -12505290
$ws.Range("N136").Value = -35362269

$ws = $wb.Worksheets.Item("CUL")
# Row 26: A Grape Idea / Grape Juice
$ws.Range("H26").Value = 269.43478
$ws.Range("I26").Value = 168
$ws.Range("K26").Value = 504
$ws.Range("M26").Value = -216
# Row 51: The Perks of Life at Sea / Jerked Beef
$ws.Range("H51").Value = 901
$ws.Range("I51").Value = 901
$ws.Range("K51").Value = 2703
$ws.Range("M51").Value = -2243
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 2624.75
$ws.Range("I113").Value = 1503
$ws.Range("K113").Value = 4509
$ws.Range("M113").Value = -2339
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 3368559.8
$ws.Range("I131").Value = 1096.5834
$ws.Range("J131").Value = 5292824.5
$ws.Range("K131").Value = 3289.7502
$ws.Range("L131").Value = 15878473.5
$ws.Range("M131").Value = 1750.2498
$ws.Range("N131").Value = -15888553.5

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 99.13333
$ws.Range("I2").Value = 52.846153
$ws.Range("K2").Value = 52.846153
$ws.Range("M2").Value = 60.153847
# Row 49: Faith and Fashion / Mythril Earrings
$ws.Range("H49").Value = 24784.25
$ws.Range("J49").Value = 25035.285
$ws.Range("L49").Value = 25035.285
$ws.Range("N49").Value = -25403.285
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 2857
$ws.Range("I113").Value = 2927.3333
$ws.Range("K113").Value = 2927.3333
$ws.Range("M113").Value = -757.3332999999998
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 6374.5
$ws.Range("I126").Value = 7166
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 21498
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -19028
$ws.Range("N126").Value = -16940
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 27751.6
$ws.Range("I132").Value = 29630.285
$ws.Range("K132").Value = 88890.855
$ws.Range("M132").Value = -86360.855

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 2840.611
$ws.Range("I93").Value = 1424.8334
$ws.Range("J93").Value = 3548.5
$ws.Range("K93").Value = 1424.8334
$ws.Range("L93").Value = 3548.5
$ws.Range("M93").Value = -176.8334
$ws.Range("N93").Value = -6044.5
# Row 114: A Heady Endeavor / Atrociraptorskin Headgear of Scouting
$ws.Range("H114").Value = 61851.668
$ws.Range("J114").Value = 61851.668
$ws.Range("L114").Value = 61851.668
$ws.Range("N114").Value = -70529.66800000001
# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 89311.336
$ws.Range("J133").Value = 89311.336
$ws.Range("L133").Value = 89311.336
$ws.Range("N133").Value = -94371.336
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 6251156
$ws.Range("I136").Value = 3907447.8
$ws.Range("K136").Value = 11722343.4
$ws.Range("M136").Value = -11719793.4

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 5377893.5
$ws.Range("I132").Value = 6667760.5
$ws.Range("K132").Value = 20003281.5
$ws.Range("M132").Value = -20000751.5
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 3879380.2
$ws.Range("I136").Value = 1740737
$ws.Range("J136").Value = 8334887
$ws.Range("K136").Value = 5222211
$ws.Range("L136").Value = 25004661
$ws.Range("M136").Value = -5219661
$ws.Range("N136").Value = -25009761

# Row 114 (BSM, Halfhearted Effort / Bluespirit Halfheart Saw) lost its HQ-profit
# figure this pull (no HQ listings on the board) - clear N114 instead of zeroing it
# so the cell is removed rather than holding a stale value.
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N114").ClearContents()
